$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.024.88'
$ws.Range('E2').Value = '  +1.35%  '
$ws.Range('D3').Value = '3.135.33'
$ws.Range('E3').Value = '  +1.28%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '534.70'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.78'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.61%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.509'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +11.32%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.32'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.79%  '
$ws.Range('E10').Value = '  +1.98%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.420'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.11%  '
$ws.Range('E12').Value = '  +3.74%  '
$ws.Range('D13').Value = '3.676.69'
$ws.Range('E13').Value = '  +1.35%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.67'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.92%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000169'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.49%  '
$ws.Range('D16').Value = '58.092.97'
$ws.Range('E16').Value = '  +1.35%  '
$ws.Range('E17').Value = '  +5.90%  '
$ws.Range('D18').Value = '3.134.03'
$ws.Range('E18').Value = '  +1.49%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.96'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.69%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.18'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.97%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '375.86'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +7.40%  '
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('E23').Value = '  -0.43%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.07'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.19%  '
$ws.Range('E25').Value = '  +2.95%  '
$ws.Range('E26').Value = '  +0.34%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.997'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.28%  '
$ws.Range('D28').Value = '0.0₃0882'
$ws.Range('E28').Value = '  +1.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.92'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +8.87%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.17'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.21%  '
$ws.Range('E31').Value = '  +0.48%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.78'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.17'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.65%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.17'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.72%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '161.22'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E36').Value = '  +3.44%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.35'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +8.83%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '25.60'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.48%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.67'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.91%  '
$ws.Range('D40').Value = '2.630.92'
$ws.Range('E40').Value = '  +9.59%  '
$ws.Range('B41').Value = 'Hedera'
$ws.Range('C41').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0674'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.33%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.23'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.47%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '38.84'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.87%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0271'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.27%  '
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('E47').Value = '  +4.08%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.977'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.40%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0996'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +10.26%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '20.22'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.40%  '
$ws.Range('E51').Value = '  -1.99%  '
